$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# Helper: split a Range [start,start+len) into several runs that preserve the
# paragraph's run-level formatting, by briefly nudging Font.Size to a
# different value and back again. Word (and this COM host) always splits a
# run at a formatting boundary, so toggling the size away and then back to
# its original value leaves the text/formatting unchanged but forces the
# desired run boundaries to persist once the document is serialized.
# ---------------------------------------------------------------------------
function Split-Runs($doc, $rangeStart, $boundaries) {
    for ($i = 0; $i -lt $boundaries.Length - 1; $i++) {
        $p1 = $rangeStart + $boundaries[$i]
        $p2 = $rangeStart + $boundaries[$i + 1]
        $r = $doc.Range($p1, $p2)
        $r.Font.Size = 12
    }
    for ($i = 0; $i -lt $boundaries.Length - 1; $i++) {
        $p1 = $rangeStart + $boundaries[$i]
        $p2 = $rangeStart + $boundaries[$i + 1]
        $r = $doc.Range($p1, $p2)
        $r.Font.Size = 11
    }
}

# ---------------------------------------------------------------------------
# 1) "2020 - Present" -> "August,2020 - February,2024" split across 5 runs:
#    "August," / "2020 - " / "Feb" / "ruary" / ",2024"
# ---------------------------------------------------------------------------
$rng = $d.Content
$rng.Find.Execute("2020 – Present") | Out-Null
$s = $rng.Start
$rng.Text = "August,2020 – February,2024"
Split-Runs $d $s @(0, 7, 14, 17, 22, 27)

# ---------------------------------------------------------------------------
# 2) "2019 - 2020 " -> "June,2019 - June,2020 " split across 7 runs:
#    "June," / "2019 " / "-" / " " / "June" / "," / "2020 "
# ---------------------------------------------------------------------------
$rng2 = $d.Content
$rng2.Find.Execute("2019 - 2020 ") | Out-Null
$s2 = $rng2.Start
$rng2.Text = "June,2019 – June,2020 "
Split-Runs $d $s2 @(0, 5, 10, 11, 12, 16, 17, 22)

# ---------------------------------------------------------------------------
# 3) Merge the "The aim of creating clusters..." run with the following
#    standalone trailing-space run into a single run (trailing space moves
#    into the sentence run).
# ---------------------------------------------------------------------------
$rng3 = $d.Content
$rng3.Find.Execute("The aim of creating clusters") | Out-Null
$s3 = $rng3.Start
$sentence = "The aim of creating clusters on a global development measurement dataset is to identify patterns, group similar entities, and uncover insights within the data. Clustering helps organize the information into meaningful segments, allowing for a more detailed analysis of different aspects of global development."
$e1 = $s3 + $sentence.Length
$e2 = $e1 + 1
$spaceRng = $d.Range($e1, $e2)
$spaceRng.Delete()
$insPoint = $d.Range($e1, $e1)
$insPoint.InsertAfter(" ")
